$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" property value (row 8, column B)
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new row for the "Jurisdiction" property right after "Contact" (row 10)
# and before "Description" (row 11), pushing the remaining property rows down.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
